# Daily attendance processing - 2026-01-18 20:35:33
# Normalizes specific "Recorded By" (column G) cell values on the active
# sheet: for the known duplicated/out-of-order author lists, the last
# author in the comma-separated list is moved to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162, column G = 7

# Exact-match replacements (old value -> new value) observed for the
# "Recorded By" column during today's processing pass.
$replacements = @{
    "dnasr281@gmail.com, System" = "System, dnasr281@gmail.com";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System";
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($null -eq $val) { continue }

    if ($replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}
